# The commit replaces the deck's "Integral" theme colour scheme with the
# default "Office Theme" colour scheme (theme1.xml, the theme driving the
# slide master / all slides). Walk the 12 theme colour slots in their
# canonical order (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) and set
# each to the standard Office palette value.
#
# COM's .RGB is a 0x00BBGGRR long, so build it from the RRGGBB hex text.
function ToComRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$officeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$tcs = $theme.ThemeColorScheme

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = ToComRgb($officeColors[$i - 1])
}
